$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data (price + volume change refresh, plus a couple of
# coin re-ordering swaps) as produced by the scheduled GitHub Actions job.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.067.25'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -5.70%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.350.82'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.98%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.63%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.00'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.07%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.350.23'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.08%  '

# Row 9
$ws.Range('E9').Value = '  -1.63%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.44'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.68%  '

# Row 11
$ws.Range('E11').Value = '  -5.16%  '

# Row 12
$ws.Range('E12').Value = '  -1.49%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.923.69'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.01%  '

# Row 14
$ws.Range('E14').Value = '  -0.29%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.347.86'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.19%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000169'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.36%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.67'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.63%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.161.88'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.54%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.70'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.67%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.47'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.84%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.10'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -7.48%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '354.69'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -7.50%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.560'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.70%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.482.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.16%  '

# Row 25
$ws.Range('E25').Value = '  +0.04%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.54%  '

# Row 27
$ws.Range('E27').Value = '  +2.27%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +18.93%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.59'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.17%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.11%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.97'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.56%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.155'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.92%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.54%  '

# Row 34
$ws.Range('E34').Value = '  -0.09%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.380.17'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.99%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.98'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.72%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.44'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.72%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.89'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.19%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.52'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.19%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '158.10'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.53%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0769'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.36%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()

# Row 43
$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.20'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +9.31%  '

# Row 44
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.18%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.751'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.59%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.75'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.63%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.70'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.67%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.58'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.18%  '

# Row 49
$ws.Range('E49').Value = '  +1.69%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.59'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +11.55%  '

# Row 51
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.893'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.36%  '
